# AddVerifyRemoveApplicant.xlsx - ValidationTestData sheet gains a block of
# "Liability" related columns (AR:BA) on the header/sample rows, mirroring
# the existing PrimaryDeclaration-style columns already on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ValidationTestData")

# --- Row 1 (field-name header row) ----------------------------------------
# Written in this exact order so the shared-string table grows the same way
# the original authoring session did (LiabilityAccountNumber was typed in
# after the Category/Type/Name/Value block, not in left-to-right order).
$ws.Range("AR1").Value = "PercentageInclude"
$ws.Range("AS1").Value = "Frequency"
$ws.Range("AT1").Value = "Balance"
$ws.Range("AU1").Value = "Limit"
$ws.Range("AV1").Value = "Payment"
$ws.Range("AX1").Value = "Category"
$ws.Range("AY1").Value = "Type"
$ws.Range("AZ1").Value = "Name"
$ws.Range("BA1").Value = "Value"
$ws.Range("AW1").Value = "LiabilityAccountNumber"

# --- Row 3 (sample data row) -----------------------------------------------
$ws.Range("AR3").Value = "100"
$ws.Range("AS3").Value = "Monthly"
$ws.Range("AT3").Value = "15000.00"
$ws.Range("AU3").Value = "10000.00"
$ws.Range("AV3").Value = "550.00"
$ws.Range("AW3").Value = "100100-2"
$ws.Range("AX3").Value = "Medical"
$ws.Range("AY3").Value = "Collection"
$ws.Range("AZ3").Value = "MICHAEL CHACOMMON"
$ws.Range("BA3").Value = "100"

# --- Formatting --------------------------------------------------------
# Row 1 headers use the same bold/fill/border style as the other headers
# (e.g. Q1). Row 2 (blank spacer row) uses the same style as the rest of
# that row (e.g. C2) - this also normalizes the old Y2/AH2 formatting so
# the whole spacer row is consistent again.
$ws.Range("Q1").Copy() | Out-Null
$ws.Range("R1:BA1").PasteSpecial(-4122) | Out-Null

$ws.Range("C2").Copy() | Out-Null
$ws.Range("R2:BA2").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- View state --------------------------------------------------------
$ws.Activate() | Out-Null
$ws.Range("L10").Select() | Out-Null
